$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B34").Value = "c61e0c5fa0c3d3aeb7f195c62229f494"
$ws.Range("B44").Value = "a2cfcbfef9b7b4aed5ed06cdf76e820f"
$ws.Range("B74").Value = "9555bf74da8a390313ded720eb47dce7"
$ws.Range("B89").Value = "160ee88f449d69ffbf488ebe9d2dcc44"
$ws.Range("B99").Value = "ec5bd2a050b8a245967e920be6cdaaa2"
$ws.Range("B110").Value = "4050bd447a74401c61ea746f9711d4fc"
$ws.Range("B161").Value = "9bb4c7968671c6ffbee5b3db18131f17"
$ws.Range("B162").Value = "28b7081ddd8b2bf574091a34d8703cef"
$ws.Range("B168").Value = "36c8cd53ba8a46717318adc0a51706b1"
$ws.Range("B180").Value = "4452182d4a3e39871668d09fdb6c1e5b"
$ws.Range("B213").Value = "e11742ebab986b101aaf472dd8371e81"
$ws.Range("B278").Value = "4f4e6e1d7f91885a3a4f184b8ac396e3"
$ws.Range("B293").Value = "21201fdc44ce87e98d9209da669acf6b"
$ws.Range("B345").Value = "183913fecc02620ae6913e0667b17656"
$ws.Range("B461").Value = "b11b80ec3b93464d6b97a5f9c1948435"
$ws.Range("B506").Value = "51d94fbb108c060af0774f3dfc25fd2e"
$ws.Range("B514").Value = "1ff4dd27e25e4cecffa8c888a063c5c2"
$ws.Range("B524").Value = "586802b4d9ba45de50d961c63708f3c0"
$ws.Range("B547").Value = "12134a6651c6de21c72dc6c1e1dae89a"
$ws.Range("B572").Value = "f1eff8d1240251c266d684e4cbc1fca7"
$ws.Range("B666").Value = "6a504f8d367e29df8fe91b6e061f2350"
$ws.Range("B729").Value = "27ed38bf1fbffac7273df8279ccba7ca"
$ws.Range("B768").Value = "8a866f38cea4d509d812189b47eef642"
$ws.Range("B816").Value = "1951623ae9020a139ec3467817acc2ab"
$ws.Range("B825").Value = "76fb08e3968f1341beee8c4d704ab1a6"
$ws.Range("B827").Value = "fe391b223dd9b3e7fc6a5f6ebd9890a3"
